# Apply updated np (non-parametric) GLMM multi-comparison results to the
# habitat*season worksheet, per the refreshed "updated criteria np" run.
# Only the ratio/std_error/statistic/adj_p_value numbers change for most rows;
# rows 7-11 additionally have their hab_1/hab_2 text re-paired.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.937753785958414
$ws.Range("H2").Value = 0.0118179857101166
$ws.Range("K2").Value = -5.09963573398489
$ws.Range("L2").Value = 0.0000330098526823834
# Row 3
$ws.Range("G3").Value = 1.04224748581956
$ws.Range("H3").Value = 0.0111770234331982
$ws.Range("K3").Value = 3.85859459487315
$ws.Range("L3").Value = 0.00912327419939958
# Row 4
$ws.Range("G4").Value = 0.829302359624809
$ws.Range("H4").Value = 0.0415215512289751
$ws.Range("K4").Value = -3.73832145877925
$ws.Range("L4").Value = 0.0144497326185885
# Row 5
$ws.Range("G5").Value = 1.23795953810776
$ws.Range("H5").Value = 0.0847759793025511
$ws.Range("K5").Value = 3.11716130182433
$ws.Range("L5").Value = 0.116865045204198
# Row 6
$ws.Range("G6").Value = 0.864338299303187
$ws.Range("H6").Value = 0.042830361884512
$ws.Range("K6").Value = -2.94213663311706
$ws.Range("L6").Value = 0.202092765590335
# Row 7
$ws.Range("D7").Value = "Shallow/Low SAV"
$ws.Range("G7").Value = 0.977371525732932
$ws.Range("H7").Value = 0.00878259303611696
$ws.Range("K7").Value = -2.54714034420563
$ws.Range("L7").Value = 0.629936513244212
# Row 8
$ws.Range("D8").Value = "Exposed/Low SAV"
$ws.Range("G8").Value = 0.841907553305544
$ws.Range("H8").Value = 0.0573066367318588
$ws.Range("K8").Value = -2.52814899368957
$ws.Range("L8").Value = 0.653594372101185
# Row 9
$ws.Range("B9").Value = "Shallow/Dense SAV"
$ws.Range("D9").Value = "Shallow/Low SAV"
$ws.Range("G9").Value = 1.13077428886452
$ws.Range("H9").Value = 0.0565624507389612
$ws.Range("K9").Value = 2.4570206735131
$ws.Range("L9").Value = 0.784529896527367
# Row 10
$ws.Range("D10").Value = "Shallow/Dense SAV"
$ws.Range("G10").Value = 1.0266427660728
$ws.Range("H10").Value = 0.086160769897687
$ws.Range("K10").Value = 0.313304696823504
# Row 11
$ws.Range("B11").Value = "Exposed/Low SAV"
$ws.Range("G11").Value = 1.16090124372388
$ws.Range("H11").Value = 0.0793471022436747
$ws.Range("K11").Value = 2.18284672501589
# Row 12
$ws.Range("G12").Value = 0.744691004084207
$ws.Range("H12").Value = 0.0110447884547377
$ws.Range("K12").Value = -19.8758367731994
$ws.Range("L12").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000121217142922633
# Row 13
$ws.Range("G13").Value = 1.48278350641426
$ws.Range("H13").Value = 0.036132355290167
$ws.Range("K13").Value = 16.1655574123271
$ws.Range("L13").Value = 0.0000000000000000000000000000000000000000000000000000000156178087235447
# Row 14
$ws.Range("G14").Value = 0.851478600656257
$ws.Range("H14").Value = 0.00972293892610806
$ws.Range("K14").Value = -14.0802596584114
$ws.Range("L14").Value = 0.000000000000000000000000000000000000000000843778509896012
# Row 15
$ws.Range("G15").Value = 0.771116300373981
$ws.Range("H15").Value = 0.0174390091270184
$ws.Range("K15").Value = -11.4929420216095
$ws.Range("L15").Value = 0.000000000000000000000000000214728179759406
# Row 16
$ws.Range("G16").Value = 0.859469220379683
$ws.Range("H16").Value = 0.0125064298533088
$ws.Range("K16").Value = -10.4073063835237
$ws.Range("L16").Value = 0.0000000000000000000000321477841792895
# Row 17
$ws.Range("G17").Value = 1.2847644942344
$ws.Range("H17").Value = 0.0310480636061752
$ws.Range("K17").Value = 10.3687759034605
$ws.Range("L17").Value = 0.0000000000000000000000478023409856087
# Row 18
$ws.Range("G18").Value = 1.14339853172172
$ws.Range("H18").Value = 0.0151058661032615
$ws.Range("K18").Value = 10.1431532819539
$ws.Range("L18").Value = 0.000000000000000000000483394671503195
# Row 19
$ws.Range("G19").Value = 1.1541286461982
$ws.Range("H19").Value = 0.017770959113091
$ws.Range("K19").Value = 9.30953185014183
$ws.Range("L19").Value = 0.00000000000000000165632157471205
# Row 20
$ws.Range("G20").Value = 1.10421553823113
$ws.Range("H20").Value = 0.0255819596676837
$ws.Range("K20").Value = 4.27905400778957
$ws.Range("L20").Value = 0.00170797281807653
# Row 21
$ws.Range("G21").Value = 0.990702843645877
$ws.Range("H21").Value = 0.0127719719361084
$ws.Range("K21").Value = -0.724539897175246
# Row 22
$ws.Range("G22").Value = 1.47635304284155
$ws.Range("H22").Value = 0.0514087535023429
$ws.Range("K22").Value = 11.1877847624563
$ws.Range("L22").Value = 0.0000000000000000000000000068792974117305
# Row 23
$ws.Range("G23").Value = 0.725718872460619
$ws.Range("H23").Value = 0.0243419256527408
$ws.Range("K23").Value = -9.55799798592138
$ws.Range("L23").Value = 0.000000000000000000158474665737375
# Row 24
$ws.Range("G24").Value = 1.35362656765584
$ws.Range("H24").Value = 0.0447220510891911
$ws.Range("K24").Value = 9.16462849629289
$ws.Range("L24").Value = 0.0000000000000000063638946489731
# Row 25
$ws.Range("G25").Value = 1.32171445651235
$ws.Range("H25").Value = 0.0454382526534183
$ws.Range("K25").Value = 8.11354811738862
$ws.Range("L25").Value = 0.0000000000000585036305152175
# Row 26
$ws.Range("G26").Value = 1.11699848296829
$ws.Range("H26").Value = 0.0161001595803721
$ws.Range("K26").Value = 7.67635112157279
$ws.Range("L26").Value = 0.00000000000189873540077491
# Row 27
$ws.Range("G27").Value = 0.9168718649101
$ws.Range("H27").Value = 0.0116492671697121
$ws.Range("K27").Value = -6.83073545492741
$ws.Range("L27").Value = 0.000000000887044164979536
# Row 28
$ws.Range("G28").Value = 1.07141726560478
$ws.Range("H28").Value = 0.0145591916633741
$ws.Range("K28").Value = 5.07643897130334
$ws.Range("L28").Value = 0.000036919127328089
# Row 29
$ws.Range("G29").Value = 0.959193125095044
$ws.Range("H29").Value = 0.0118101094290708
$ws.Range("K29").Value = -3.38377154414893
$ws.Range("L29").Value = 0.0493332527255936
# Row 30
$ws.Range("G30").Value = 0.982352346411932
$ws.Range("H30").Value = 0.00848986886056728
$ws.Range("K30").Value = -2.06022140475075
# Row 31
$ws.Range("G31").Value = 1.02414448218089
$ws.Range("H31").Value = 0.0116939121348173
$ws.Range("K31").Value = 2.08943268627799
# Row 32
$ws.Range("G32").Value = 1.12465028522679
$ws.Range("H32").Value = 0.0175786935438224
$ws.Range("K32").Value = 7.51563612532647
$ws.Range("L32").Value = 0.00000000000645637549363058
# Row 33
$ws.Range("G33").Value = 0.896352630794693
$ws.Range("H33").Value = 0.0130612107098371
$ws.Range("K33").Value = -7.50926893749397
$ws.Range("L33").Value = 0.00000000000671880650283911
# Row 34
$ws.Range("G34").Value = 0.890148522155623
$ws.Range("H34").Value = 0.0138583328828509
$ws.Range("K34").Value = -7.47448272815231
$ws.Range("L34").Value = 0.00000000000868091004886681
# Row 35
$ws.Range("G35").Value = 1.12664712780148
$ws.Range("H35").Value = 0.0333002800967375
$ws.Range("K35").Value = 4.03444810194608
$ws.Range("L35").Value = 0.00459738591154588
# Row 36
$ws.Range("G36").Value = 1.12540266953022
$ws.Range("H36").Value = 0.0333318443299396
$ws.Range("K36").Value = 3.98886071295817
$ws.Range("L36").Value = 0.00544409387893987
# Row 37
$ws.Range("G37").Value = 0.894763956709487
$ws.Range("H37").Value = 0.0260335386423453
$ws.Range("K37").Value = -3.82174608340985
$ws.Range("L37").Value = 0.0104682973231737
# Row 38
$ws.Range("G38").Value = 1.00696974548028
$ws.Range("H38").Value = 0.00748040531428786
$ws.Range("K38").Value = 0.934973123735247
# Row 39
$ws.Range("G39").Value = 1.00808324188703
$ws.Range("H39").Value = 0.00638072911190895
$ws.Range("K39").Value = 1.27192731011726
# Row 40
$ws.Range("G40").Value = 1.00177552311232
$ws.Range("H40").Value = 0.0323473072425048
$ws.Range("K40").Value = 0.0549380635642719
# Row 41
$ws.Range("G41").Value = 1.00110578933653
$ws.Range("H41").Value = 0.00932674688805125
$ws.Range("K41").Value = 0.118626624010014
